$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 13-14; existing rows 13:39 shift down to 15:41.
$ws.Rows("13:14").Insert()

# New row 13 - "Primera" quality entry for date 44725 (2022-06-13)
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C13").Value = "Arica y Parinacota"
$ws.Range("D13").Value = 44725
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100108
$ws.Range("H13").Value = "Tropicales y subtropicales"
$ws.Range("I13").Value = 100108001
$ws.Range("J13").Value = "Guayaba"
$ws.Range("K13").Value = "Sin especificar"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 140
$ws.Range("N13").Value = 700
$ws.Range("O13").Value = 800
$ws.Range("P13").Value = 750
$ws.Range("Q13").Value = "$/kilo (en caja de 10 kilos )"
$ws.Range("R13").Value = "Región de Arica y Parinacota"
$ws.Range("S13").Value = 750
$ws.Range("T13").Value = 1

# New row 14 - "Segunda" quality entry for date 44725 (2022-06-13)
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C14").Value = "Arica y Parinacota"
$ws.Range("D14").Value = 44725
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100108
$ws.Range("H14").Value = "Tropicales y subtropicales"
$ws.Range("I14").Value = 100108001
$ws.Range("J14").Value = "Guayaba"
$ws.Range("K14").Value = "Sin especificar"
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 500
$ws.Range("O14").Value = 600
$ws.Range("P14").Value = 550
$ws.Range("Q14").Value = "$/kilo (en caja de 10 kilos )"
$ws.Range("R14").Value = "Región de Arica y Parinacota"
$ws.Range("S14").Value = 550
$ws.Range("T14").Value = 1
